$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 34; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = $cell.Value2 + 1
}
